# Rename sheets: Sheet1 -> Old, Sheet2 -> Current
$wb = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item(1)
$wsOld.Name = "Old"
$wsCurrent = $wb.Worksheets.Item(2)
$wsCurrent.Name = "Current"

# Fix up the Print_Area defined name so it points at the renamed sheet
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Old!Print_Area") {
        $n.RefersTo = "=Old!`$A`$1:`$E`$55"
    }
}
